$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("B78:AP78").ClearContents()
